# Actualización 11 de Mayo - Tarde
# Updates exam-result stats for three groups (rows 4, 6, 7) on the
# "2o Parcial" and "3er Parcial" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "2o Parcial": rows 4, 6, 7 go from "no grades submitted yet"
# (Aprobados=0, Reprobados=Totales, no Promedio) to graded results,
# and gain a Promedio (column I) value.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 4 - 4APM
$ws2.Range("E4").Value = 16
$ws2.Range("F4").Value = 12
$ws2.Range("G4").Value = 57.14
$ws2.Range("H4").Value = 42.86
$ws2.Range("I4").Value = 6.2
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0

# Row 6 - 4BLCM
$ws2.Range("E6").Value = 28
$ws2.Range("F6").Value = 8
$ws2.Range("G6").Value = 77.78
$ws2.Range("H6").Value = 22.22
$ws2.Range("I6").Value = 8
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0

# Row 7 - 4BEM
$ws2.Range("E7").Value = 21
$ws2.Range("F7").Value = 16
$ws2.Range("G7").Value = 56.76
$ws2.Range("H7").Value = 43.24
$ws2.Range("I7").Value = 6.4
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0

# ---------------------------------------------------------------
# Sheet "3er Parcial": rows 4, 6, 7 get refreshed Aprobados/
# Reprobados/percentage/Promedio figures.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3er Parcial")

# Row 4 - 4APM
$ws3.Range("E4").Value = 16
$ws3.Range("F4").Value = 12
$ws3.Range("G4").Value = 57.14
$ws3.Range("H4").Value = 42.86

# Row 6 - 4BLCM (only Promedio changes)
$ws3.Range("I6").Value = 8.300000000000001

# Row 7 - 4BEM
$ws3.Range("E7").Value = 21
$ws3.Range("F7").Value = 16
$ws3.Range("G7").Value = 56.76
$ws3.Range("H7").Value = 43.24
$ws3.Range("I7").Value = 6.8
